$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.391.68"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.794.31"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.04%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "307.11"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4549"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.61%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3599"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.82%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.35"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07117"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8837"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07820"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.50"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.808.15"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.49%  "
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.284"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.324"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "85.07"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.008"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008579"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.006"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.28"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.09%  "
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "WrappedBTC"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "26.420.28"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.988"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.039.29"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.61%  "
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.53"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.977"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "152.51"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.90"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.85%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.047"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "111.91"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.871"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08654"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.16%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.449"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7261"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.712"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.110"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.48%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.076"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01943"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.05112"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.877"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5173"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.891"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1525"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.008"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4671"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.006"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.904"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.70%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "100.18"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.12%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.588"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.34%  "
